$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.32"
$ws.Range("E2").Value = "'-0.65%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'36.74"
$ws.Range("E3").Value = "'2.84%"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'5.005"
$ws.Range("E4").Value = "'-0.85%"
$ws.Range("G4").Value = "'21"
$ws.Range("D5").Value = "'0.07712"
$ws.Range("E5").Value = "'-1.06%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'2.058"
$ws.Range("E6").Value = "'-8.70%"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'7.973"
$ws.Range("E7").Value = "'-1.70%"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'4.031"
$ws.Range("E8").Value = "'-0.22%"
$ws.Range("G8").Value = "'21"
$ws.Range("D9").Value = "'0.9166"
$ws.Range("E9").Value = "'-1.30%"
$ws.Range("G9").Value = "'21"
$ws.Range("D10").Value = "'0.09672"
$ws.Range("E10").Value = "'1.39%"
$ws.Range("G10").Value = "'21"
$ws.Range("D11").Value = "'0.1849"
$ws.Range("E11").Value = "'0.98%"
$ws.Range("G11").Value = "'21"
$ws.Range("D12").Value = "'0.08575"
$ws.Range("E12").Value = "'0.28%"
$ws.Range("G12").Value = "'21"
$ws.Range("D13").Value = "'0.03577"
$ws.Range("E13").Value = "'4.36%"
$ws.Range("G13").Value = "'21"
$ws.Range("D14").Value = "'0.09935"
$ws.Range("E14").Value = "'-0.18%"
$ws.Range("G14").Value = "'21"
$ws.Range("D15").Value = "'0.001471"
$ws.Range("E15").Value = "'-1.03%"
$ws.Range("G15").Value = "'21"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005743"
$ws.Range("E16").Value = "'0.02%"
$ws.Range("G16").Value = "'21"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.472"
$ws.Range("E17").Value = "'-0.24%"
$ws.Range("G17").Value = "'21"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.377"
$ws.Range("E18").Value = "'9.10%"
$ws.Range("G18").Value = "'21"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3434"
$ws.Range("E19").Value = "'0.82%"
$ws.Range("G19").Value = "'21"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1338"
$ws.Range("E20").Value = "'1.15%"
$ws.Range("G20").Value = "'21"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'4.966"
$ws.Range("E21").Value = "'9.06%"
$ws.Range("G21").Value = "'21"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2212"
$ws.Range("E22").Value = "'-1.22%"
$ws.Range("G22").Value = "'21"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04612"
$ws.Range("E23").Value = "'-1.48%"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'0.005107"
$ws.Range("E24").Value = "'12.53%"
$ws.Range("G24").Value = "'21"
$ws.Range("D25").Value = "'0.001238"
$ws.Range("E25").Value = "'-0.48%"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.0001407"
$ws.Range("E26").Value = "'7.95%"
$ws.Range("G26").Value = "'21"
$ws.Range("G27").Value = "'21"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("D39").Value = "'0.01751"
$ws.Range("E39").Value = "'-1.03%"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.04604"
$ws.Range("E40").Value = "'-2.39%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.007696"
$ws.Range("E41").Value = "'-3.57%"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.1386"
$ws.Range("E42").Value = "'-2.48%"
$ws.Range("G42").Value = "'21"
$ws.Range("D43").Value = "'0.007976"
$ws.Range("E43").Value = "'-0.52%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.002163"
$ws.Range("E44").Value = "'-2.55%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.009862"
$ws.Range("E45").Value = "'8.13%"
$ws.Range("G45").Value = "'21"
$ws.Range("D46").Value = "'0.00006289"
$ws.Range("E46").Value = "'1.13%"
$ws.Range("G46").Value = "'21"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.25%"
$ws.Range("G47").Value = "'21"
$ws.Range("D48").Value = "'0.0005853"
$ws.Range("E48").Value = "'0.90%"
$ws.Range("G48").Value = "'21"
$ws.Range("D49").Value = "'34.05"
$ws.Range("E49").Value = "'553.40%"
$ws.Range("G49").Value = "'21"
$ws.Range("D50").Value = "'0.002010"
$ws.Range("E50").Value = "'-25.47%"
$ws.Range("G50").Value = "'21"
$ws.Range("D51").Value = "'0.00002111"
$ws.Range("E51").Value = "'0.25%"
$ws.Range("G51").Value = "'21"
